$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (prevents Excel's automatic date /
# number inference for strings like "02-18-2022"), then restore the cell's
# style to Normal so no stray NumberFormat/style is left behind.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 2 — update in place (fuzzy-matched record AN22-000188 -> AN22-000195)
# ---------------------------------------------------------------------
Set-TextValue "A2" "AN22-000195"
$ws.Range("B2").Value = "Removing CMU Paint from Classroom"
$ws.Range("C2").Value = "Notification"
$ws.Range("D2").Value = "Christine Del Viscio"
$ws.Range("F2").Value = "2501 S 63RD ST"
Set-TextValue "G2" "02-18-2022"
Set-TextValue "H2" "02-18-2022"
Set-TextValue "K2" "02-19-2022"
Set-TextValue "L2" "03-04-2022"
$ws.Range("Q2").Value = "Thomas Morton Sch"
$ws.Range("R2").Value = "School District of Philadelphia"
$ws.Range("S2").Value = "440 N Board Street Philadelphia Pa"
$ws.Range("T2").Value = "DELTA / B.J.D.S., INC."
$ws.Range("Y2").Value = 50

# ---------------------------------------------------------------------
# Row 3 — new record AN22-000197
# ---------------------------------------------------------------------
Set-TextValue "A3" "AN22-000197"
$ws.Range("B3").Value = "GIRLS HS - CAP 17SF OF EXTERIOR BOILER INSULATION, 60 SF OF BOILER BREEHING & REM 6 SF OF MATERIAL ASSOC W/ EXT. BOILER INSULATION (DEBRIS) IN MAIN MECHANICAL ROOM LOWER LEVEL, CAP <1 LF ACPFI IN 1ST ROOM IN MECHANICAL SPACE AT ENTRANCE, <1 LF ACPFI IN PUMP UTILITY ROOM, <1 LF ACPFI * CAP 2 SF OF TANK INSULATION IN MAIN MECANICAL ROOM UPPER LEVEL, REM 8 LF ACPFI IN PIPE CHASE IN HALL FROM PHY ED DEPT OFC TO VISITR CHANGE RM, REM 4 LF ACPFI IN PIPE CHASE IN GIRLS LOCKER ROOM, 8 LF ACPFI IN PC IN"
$ws.Range("C3").Value = "Notification with Alternative Methods"
$ws.Range("D3").Value = "Stephen Link"
$ws.Range("E3").Value = "Issued-ASB"
$ws.Range("F3").Value = "1400 W OLNEY AVE"
Set-TextValue "G3" "02-18-2022"
Set-TextValue "H3" "02-23-2022"
Set-TextValue "I3" "02-17-2022"
Set-TextValue "J3" "04-04-2022"
Set-TextValue "K3" "02-17-2022"
Set-TextValue "L3" "04-04-2022"
$ws.Range("N3").Value = 4483018
$ws.Range("O3").Value = "Major Removal Project"
$ws.Range("P3").Value = "Renovation"
$ws.Range("Q3").Value = "GIRLS HS"
$ws.Range("R3").Value = "SCHOOL DIST OF PHILA"
$ws.Range("S3").Value = "440 N BROAD STREET SUITE 373 PHILA., PA 19130"
$ws.Range("T3").Value = "SCHOOL DISTRICT OF PHILADELPHIA"
$ws.Range("W3").Value = "Yes"
$ws.Range("X3").Value = 69
$ws.Range("Y3").Value = 87
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = "Yes"

# ---------------------------------------------------------------------
# Row 4 — new record AN22-000193
# ---------------------------------------------------------------------
Set-TextValue "A4" "AN22-000193"
$ws.Range("B4").Value = "Asbestos Abatement"
$ws.Range("C4").Value = "Notification with Alternative Methods"
$ws.Range("D4").Value = "Chrissy Bruno"
$ws.Range("E4").Value = "Approved-ASB"
$ws.Range("F4").Value = "2901 PRINCETON AVE"
Set-TextValue "G4" "02-18-2022"
Set-TextValue "H4" "02-18-2022"
Set-TextValue "K4" "02-21-2022"
Set-TextValue "L4" "07-21-2022"
$ws.Range("O4").Value = "Minor Removal Project"
$ws.Range("P4").Value = "Renovation"
$ws.Range("Q4").Value = "MAYFAIR ELEMENTARY SCHOOL"
$ws.Range("R4").Value = "School District of Philadelphia"
$ws.Range("S4").Value = "440 NORTH BROAD STREET PHILADELPHIA, PA 19130"
$ws.Range("T4").Value = "PEPPER ENVIRONMENTAL SERVICES"
$ws.Range("W4").Value = "Yes"
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 16
$ws.Range("AB4").Value = 20
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = "No"

# ---------------------------------------------------------------------
# Row 5 — new record AN22-000194
# ---------------------------------------------------------------------
Set-TextValue "A5" "AN22-000194"
$ws.Range("B5").Value = "Asbestos Abatement"
$ws.Range("C5").Value = "Notification with Alternative Methods"
$ws.Range("D5").Value = "Chrissy Bruno"
$ws.Range("E5").Value = "Approved-ASB"
$ws.Range("F5").Value = "5925 MALVERN AVE"
Set-TextValue "G5" "02-18-2022"
Set-TextValue "H5" "02-18-2022"
Set-TextValue "K5" "02-21-2022"
Set-TextValue "L5" "02-21-2023"
$ws.Range("O5").Value = "Non Friable Removal"
$ws.Range("P5").Value = "Renovation"
$ws.Range("Q5").Value = "SCIENCE LEADERSHIP ACADEMY AT BEEBER"
$ws.Range("R5").Value = "School District of Philadelphia"
$ws.Range("S5").Value = "440 NORTH BROAD STREET PHILADELPHIA, PA 19130"
$ws.Range("T5").Value = "PEPPER ENVIRONMENTAL SERVICES"
$ws.Range("W5").Value = "Yes"
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 20
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = "No"

# ---------------------------------------------------------------------
# Row 6 — new record AN22-000204
# ---------------------------------------------------------------------
Set-TextValue "A6" "AN22-000204"
$ws.Range("B6").Value = "FITZPATRICK ES - REM 6 LF ACPFI & 70 SF OF VDC IN AUDITORIUM STAGE 115 & 2 LF ACPFI IN BE OFFICE 017B."
$ws.Range("C6").Value = "Notification"
$ws.Range("D6").Value = "Stephen Link"
$ws.Range("E6").Value = "Approved"
$ws.Range("F6").Value = "11061 KNIGHTS RD"
Set-TextValue "G6" "02-22-2022"
Set-TextValue "H6" "02-23-2022"
Set-TextValue "K6" "02-23-2022"
Set-TextValue "L6" "03-07-2022"
$ws.Range("O6").Value = "Minor Removal Project"
$ws.Range("P6").Value = "Renovation"
$ws.Range("Q6").Value = "FITZPATRICK ES"
$ws.Range("R6").Value = "SCHOOL DIST OF PHILA"
$ws.Range("S6").Value = "440 N BROAD STREET SUITE 373 PHILA., PA 19130"
$ws.Range("T6").Value = "SCHOOL DISTRICT OF PHILADELPHIA"
$ws.Range("W6").Value = "Yes"
$ws.Range("X6").Value = 8
$ws.Range("Y6").Value = 70
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = "No"
